$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update column F ("想去人数") for specific rows
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 65
$wsExpo.Range("F4").Value = 163
$wsExpo.Range("F5").Value = 359
$wsExpo.Range("F6").Value = 5396
$wsExpo.Range("F8").Value = 5387
$wsExpo.Range("F9").Value = 629
$wsExpo.Range("F11").Value = 1378
$wsExpo.Range("F12").Value = 8

# Sheet "全部类型" (all types) - update column F ("想去人数") for specific rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 65
$wsAll.Range("F4").Value = 163
$wsAll.Range("F6").Value = 359
$wsAll.Range("F7").Value = 5396
$wsAll.Range("F9").Value = 5387
$wsAll.Range("F10").Value = 629
$wsAll.Range("F12").Value = 1378
$wsAll.Range("F13").Value = 8
